$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2201492537313433
$ws.Range("C2").Value = 0.458955223880597
$ws.Range("J2").Value = 0.01865671641791045
$ws.Range("P2").Value = 0.1940298507462687
$ws.Range("S2").Value = 0.1082089552238806
$ws.Range("B3").Value = 0.01612903225806452
$ws.Range("C3").Value = 0.03225806451612903
$ws.Range("J3").Value = 0.02419354838709677
$ws.Range("O3").Value = 0.008064516129032258
$ws.Range("P3").Value = 0.7983870967741935
$ws.Range("S3").Value = 0.1209677419354839
$ws.Range("J4").Value = 0.04651162790697674
$ws.Range("O4").Value = 0.02325581395348837
$ws.Range("P4").Value = 0.6744186046511628
$ws.Range("S4").Value = 0.2558139534883721
$ws.Range("B6").Value = 0.075
$ws.Range("D6").Value = 0.01
$ws.Range("F6").Value = 0.06
$ws.Range("J6").Value = 0.255
$ws.Range("O6").Value = 0.015
$ws.Range("Q6").Value = 0.15
$ws.Range("R6").Value = 0.06
$ws.Range("S6").Value = 0.375
$ws.Range("B7").Value = 0.1169354838709677
$ws.Range("D7").Value = 0.01209677419354839
$ws.Range("E7").Value = 0.008064516129032258
$ws.Range("F7").Value = 0.0282258064516129
$ws.Range("J7").Value = 0.1330645161290323
$ws.Range("O7").Value = 0.004032258064516129
$ws.Range("Q7").Value = 0.2217741935483871
$ws.Range("R7").Value = 0.04838709677419355
$ws.Range("S7").Value = 0.4274193548387097
$ws.Range("B8").Value = 0.06691449814126393
$ws.Range("D8").Value = 0.0241635687732342
$ws.Range("F8").Value = 0.05947955390334572
$ws.Range("J8").Value = 0.1152416356877323
$ws.Range("O8").Value = 0.01115241635687732
$ws.Range("Q8").Value = 0.1802973977695167
$ws.Range("R8").Value = 0.08736059479553904
$ws.Range("S8").Value = 0.4553903345724907
$ws.Range("B9").Value = 0.06944444444444445
$ws.Range("F9").Value = 0.08333333333333333
$ws.Range("J9").Value = 0.1111111111111111
$ws.Range("O9").Value = 0.01388888888888889
$ws.Range("Q9").Value = 0.1944444444444444
$ws.Range("R9").Value = 0.0625
$ws.Range("S9").Value = 0.4652777777777778
$ws.Range("B10").Value = 0.09179528838342811
$ws.Range("D10").Value = 0.02274573517465475
$ws.Range("F10").Value = 0.06580016246953696
$ws.Range("J10").Value = 0.1121039805036556
$ws.Range("O10").Value = 0.008935824532900082
$ws.Range("Q10").Value = 0.2347684809098294
$ws.Range("R10").Value = 0.06580016246953696
$ws.Range("S10").Value = 0.3980503655564582
$ws.Range("G11").Value = 0.1440922190201729
$ws.Range("J11").Value = 0.06628242074927954
$ws.Range("K11").Value = 0.1988472622478386
$ws.Range("L11").Value = 0.5706051873198847
$ws.Range("S11").Value = 0.02017291066282421
$ws.Range("G12").Value = 0.7867298578199052
$ws.Range("J12").Value = 0.1184834123222749
$ws.Range("L12").Value = 0.04265402843601896
$ws.Range("S12").Value = 0.05213270142180094
$ws.Range("G13").Value = 0.6610169491525424
$ws.Range("J13").Value = 0.2711864406779661
$ws.Range("S13").Value = 0.06779661016949153
$ws.Range("F15").Value = 0.01005025125628141
$ws.Range("H15").Value = 0.2261306532663317
$ws.Range("I15").Value = 0.05527638190954774
$ws.Range("J15").Value = 0.3065326633165829
$ws.Range("K15").Value = 0.1005025125628141
$ws.Range("M15").Value = 0.01507537688442211
$ws.Range("O15").Value = 0.06532663316582915
$ws.Range("S15").Value = 0.221105527638191
$ws.Range("F16").Value = 0.02209944751381215
$ws.Range("H16").Value = 0.1657458563535912
$ws.Range("I16").Value = 0.04419889502762431
$ws.Range("J16").Value = 0.3646408839779006
$ws.Range("K16").Value = 0.1767955801104972
$ws.Range("M16").Value = 0.01657458563535912
$ws.Range("O16").Value = 0.02209944751381215
$ws.Range("S16").Value = 0.1878453038674033
$ws.Range("F17").Value = 0.01814516129032258
$ws.Range("H17").Value = 0.217741935483871
$ws.Range("I17").Value = 0.07258064516129033
$ws.Range("J17").Value = 0.3931451612903226
$ws.Range("K17").Value = 0.1088709677419355
$ws.Range("M17").Value = 0.01814516129032258
$ws.Range("O17").Value = 0.04435483870967742
$ws.Range("S17").Value = 0.1270161290322581
$ws.Range("F18").Value = 0.01863354037267081
$ws.Range("H18").Value = 0.2049689440993789
$ws.Range("I18").Value = 0.09937888198757763
$ws.Range("J18").Value = 0.4347826086956522
$ws.Range("K18").Value = 0.09316770186335403
$ws.Range("M18").Value = 0.03105590062111801
$ws.Range("N18").Value = 0.006211180124223602
$ws.Range("O18").Value = 0.04347826086956522
$ws.Range("S18").Value = 0.06832298136645963
$ws.Range("F19").Value = 0.01111934766493699
$ws.Range("H19").Value = 0.2446256486286138
$ws.Range("I19").Value = 0.05485544848035582
$ws.Range("J19").Value = 0.3573017049666419
$ws.Range("K19").Value = 0.1134173461823573
$ws.Range("M19").Value = 0.03039288361749444
$ws.Range("N19").Value = 0.0007412898443291327
$ws.Range("O19").Value = 0.072646404744255
$ws.Range("S19").Value = 0.1148999258710156
